$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Account Affected column must stay text even though it looks numeric
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D4").NumberFormat = "@"

# Row 3: Withdraw transaction
$ws.Range("A3").Value = "2025-10-21 20:45:54"
$ws.Range("B3").Value = "moumi"
$ws.Range("C3").Value = "Withdraw"
$ws.Range("D3").Value = "12344"
$ws.Range("E3").Value = 1

# Row 4: Deposit transaction
$ws.Range("A4").Value = "2025-10-21 20:46:47"
$ws.Range("B4").Value = "moumi"
$ws.Range("C4").Value = "Deposit"
$ws.Range("D4").Value = "12344"
$ws.Range("E4").Value = 2
